$wb = $excel.ActiveWorkbook

# --- "Active" sheet (Todo list): add the two new follow-up tasks that were
#     spun off while removing the old find-regions code. They're inserted
#     right after "make sure all TODOs are cleared from project" (row 16),
#     matching the row order in the updated workbook.
$active = $wb.Worksheets.Item("Active")

$active.Rows.Item(17).Insert()
$active.Rows.Item(18).Insert()

# New rows inherit the header's bold styling from the insert; reset it back
# to the plain look the rest of the data rows use.
$active.Range("A17:E18").Font.Bold = $false

$active.Range("A17").Value = 115
$active.Range("B17").Value = "use mouse scroll for zoom in/out"
$active.Range("C17").Value = "Todo"
$active.Range("D17").Value = "Task"
$active.Range("E17").Value = "'9/10/2018"

$active.Range("A18").Value = 116
$active.Range("B18").Value = "drag and drop image to move it"
$active.Range("C18").Value = "Todo"
$active.Range("D18").Value = "Task"
$active.Range("E18").Value = "'9/10/2018"

# --- "Inactive" sheet (Done list): the "remove old code for old
#     find-regions logic" task is now finished, so it moves here as a new
#     top row (most-recently-completed tasks are listed first).
$inactive = $wb.Worksheets.Item("Inactive")

$inactive.Rows.Item(2).Insert()
$inactive.Range("A2:F2").Font.Bold = $false

$inactive.Range("A2").Value = 114
$inactive.Range("B2").Value = "remove old code for old find-regions logic"
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "Task"
$inactive.Range("E2").Value = "'9/10/2018"
$inactive.Range("F2").Value = "'9/10/2018"

# --- "Config" sheet: bump the Max Id tracker now that ids 114-116 are used.
$config = $wb.Worksheets.Item("Config")
$config.Range("F2").Value = 116
